$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Datatype Customer: add new "String dob_s" field (row 8, previously blank) ---
$ws.Range("B8").Value = "String"
$ws.Range("C8").Value = "dob_s"
$ws.Range("B4:C4").Copy()
$ws.Range("B8:C8").PasteSpecial(-4122)

# --- Datatype Customer2: add new "String dob_s" field (row 16, previously blank) ---
$ws.Range("B16").Value = "String"
$ws.Range("C16").Value = "dob_s"
$ws.Range("B12:C12").Copy()
$ws.Range("B16:C16").PasteSpecial(-4122)

# --- Method Customer proxyCustomer(Customer c): insert new body line before "return c;" ---
# Row27 held "return c;" -> becomes the new assignment line; row28 (blank) becomes "return c;"
$ws.Range("B27").Value = 'c.dob_s = toString(c.dob, "yyyy MM dd - HH mm ss");'
$ws.Range("B27").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("B28").Value = "return c;"

# --- Method Customer2 proxyCustomer2(Customer2 c): shift body down by one row ---
# Row30 held the method signature, row31 held "return c;". Clear row30, shift
# signature to row31, and add the new assignment line + return at rows 32/33.
$ws.Range("B30").ClearContents()
$ws.Range("B29").Copy()
$ws.Range("B30").PasteSpecial(-4122)

$ws.Range("B31").Value = "Method Customer2 proxyCustomer2(Customer2 c)"

$ws.Range("B31").Copy()
$ws.Range("B32:B33").PasteSpecial(-4122)
$ws.Range("B32").Value = 'c.dob_s = toString(c.dob, "yyyy MM dd - HH mm ss");'
$ws.Range("B33").Value = "return c;"

# --- Append two new blank rows at the end of the table (43, 44) ---
$ws.Range("A42:E42").Copy()
$ws.Range("A43:E44").PasteSpecial(-4122)
$ws.Rows("43:44").RowHeight = $ws.Rows(42).RowHeight

# --- Column width tweaks ---
$ws.Columns("B").ColumnWidth = 20.666666666666668
$ws.Columns("C").ColumnWidth = 17.666666666666668

# --- Give the whole table area an explicit white fill (was "no fill") ---
$ws.Range("A1:E44").Interior.ColorIndex = 2
